$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1) "Netid: _______________" heading: remove the spell-check proofErr
#    markers around "Netid" and collapse the two runs into a single run.
#    (Deleting the whole paragraph -- including its own paragraph mark --
#    and re-inserting fresh OOXML is the only way to also drop the
#    orphan <w:proofErr/> markers, since they are not addressable through
#    the normal Range/Text object model.)
# ----------------------------------------------------------------------
$netidPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Netid:*") {
        $netidPara = $p
        break
    }
}

if ($netidPara -ne $null) {
    $prevPara = $netidPara.Previous()
    $netidPara.Range.Delete()

    $insertAt = $prevPara.Range
    $insertAt.Collapse(0)
    $insertAt.InsertParagraphAfter()

    $newNetidPara = $prevPara.Next()
    $newNetidPara.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Netid: _______________</w:t></w:r></w:p>')
}

# ----------------------------------------------------------------------
# 2) Insert a new numbered list item ("Explain what challenges ...")
#    right after the last existing "Show a screenshot ..." list item,
#    and before the "For example:" paragraph.
# ----------------------------------------------------------------------
$lastListItem = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Show a screenshot your corrected solution*") {
        $lastListItem = $p
        break
    }
}

if ($lastListItem -ne $null) {
    $r = $lastListItem.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()

    $newItem = $lastListItem.Next()
    $newItem.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Explain what challenges prevented you from finding this solution while you took the test.</w:t></w:r></w:p>')
}

# ----------------------------------------------------------------------
# 3) Mark every inline picture's run as <w:noProof/> (matches what Word
#    writes once a picture has actually been rendered/laid out).
# ----------------------------------------------------------------------
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shp = $d.InlineShapes.Item($i)
    $shp.Range.NoProofing = $true
}
